# Add a 2023 column (K) to the "Number of employees" table, matching the
# look of the existing 2022 column (J) and closing the table off with a
# right-hand border on the new last column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values -----------------------------------------------------------
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 9382
$ws.Range("K5").Value = 3250
$ws.Range("K6").Value = 6132

# xlEdge* border constants / line styles
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlLineStyleNone = -4142

# The year-header row (row 3, columns B:J) used to carry both a top and a
# bottom border, which doubled up visually with row 4's top border. Clean
# that up to a single separator line, same as the rest of the "upgraded"
# tables: drop the redundant bottom border on the numeric header cells
# (columns B through the new last column, K).
$ws.Range("B3:K3").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone

foreach ($row in 3..6) {
    $src = $ws.Range("J$row")
    $dst = $ws.Range("K$row")

    # Font
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold

    # Fill: new cells already inherit the table's white fill by default, so
    # no explicit Interior write is needed here.

    # Alignment
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment

    # Number format (literal format code so the escaped space round-trips
    # correctly instead of losing its backslash through NumberFormat's
    # string round-trip)
    if ($src.NumberFormat -ne "General") {
        $dst.NumberFormat = "#\ ##0"
    }
}

# Borders for the new column, set explicitly per row (matching the rest of
# the table row-by-row) rather than copied live from column J, since the
# bottom-border clear above happens in the same pass.
$ws.Range("K3").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("K4").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("K6").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

# ... and add the closing right-hand border down the new last column (K),
# since K is now the last (rightmost) column of the table.
$ws.Range("K3:K6").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

$ws.Range("A1:K6").Select()
